$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E (text-like numeric values, e.g. "7.00") to be stored as literal text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '43.332.21'
$ws.Range('E2').Value = '  +2.71%  '

$ws.Range('D3').Value = '2.305.07'
$ws.Range('E3').Value = '  +1.62%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '310.40'
$ws.Range('E5').Value = '  +1.38%  '

$ws.Range('D6').Value = '103.11'
$ws.Range('E6').Value = '  +5.82%  '

$ws.Range('D7').Value = '0.533'
$ws.Range('E7').Value = '  +1.37%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('E9').Value = '  +8.50%  '

$ws.Range('D10').Value = '35.77'
$ws.Range('E10').Value = '  +1.37%  '

$ws.Range('E11').Value = '  +3.04%  '

$ws.Range('E12').Value = '  -1.13%  '

$ws.Range('D13').Value = '7.00'
$ws.Range('E13').Value = '  +1.76%  '

$ws.Range('D14').Value = '2.662.60'
$ws.Range('E14').Value = '  +1.65%  '

$ws.Range('D15').Value = '14.98'
$ws.Range('E15').Value = '  +2.04%  '

$ws.Range('D16').Value = '2.297.57'
$ws.Range('E16').Value = '  +0.63%  '

$ws.Range('D17').Value = '0.808'
$ws.Range('E17').Value = '  +2.21%  '

$ws.Range('D18').Value = '43.231.62'
$ws.Range('E18').Value = '  +2.74%  '

$ws.Range('D19').Value = '12.26'
$ws.Range('E19').Value = '  -0.37%  '

$ws.Range('D20').Value = '0.0₃0934'
$ws.Range('E20').Value = '  +3.23%  '

$ws.Range('D21').Value = '6.17'
$ws.Range('E21').Value = '  +2.79%  '

$ws.Range('D22').Value = '68.14'
$ws.Range('E22').Value = '  +0.60%  '

$ws.Range('D23').Value = '241.46'
$ws.Range('E23').Value = '  +1.88%  '

$ws.Range('D24').Value = '2.01'
$ws.Range('E24').Value = '  +1.36%  '

$ws.Range('D25').Value = '2.62'
$ws.Range('E25').Value = '  +1.66%  '

$ws.Range('E26').Value = '  +0.08%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '24.93'
$ws.Range('E27').Value = '  +6.06%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '2.30'
$ws.Range('E28').Value = '  +8.11%  '

$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = '36.60'
$ws.Range('E29').Value = '  -2.25%  '

$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').Value = '9.65'
$ws.Range('E30').Value = '  +0.67%  '

$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = '171.49'
$ws.Range('E31').Value = '  +5.66%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '5.27'
$ws.Range('E32').Value = '  +0.30%  '

$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.01%  '

$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = '2.58'
$ws.Range('E34').Value = '  +8.59%  '

$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').Value = '17.79'
$ws.Range('E35').Value = '  +0.52%  '

$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '0.0742'
$ws.Range('E36').Value = '  +1.08%  '

$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '3.06'
$ws.Range('E37').Value = '  -2.38%  '

$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').Value = '1.88'
$ws.Range('E38').Value = '  +3.05%  '

$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '0.106'
$ws.Range('E39').Value = '  +1.70%  '

$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = '0.116'
$ws.Range('E40').Value = '  +1.35%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '4.35'
$ws.Range('E41').Value = '  +6.65%  '

$ws.Range('B42').Value = 'ApeXProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D42').Value = '2.30'
$ws.Range('E42').Value = '  -2.06%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.0290'
$ws.Range('E43').Value = '  +3.69%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '19.25'
$ws.Range('E44').Value = '  +1.53%  '

$ws.Range('D45').Value = '1.968.06'
$ws.Range('E45').Value = '  +0.97%  '

$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = '3.00'
$ws.Range('E46').Value = '  +3.01%  '

$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '9.98'
$ws.Range('E47').Value = '  +0.36%  '

$ws.Range('B48').Value = 'MultiversX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D48').Value = '55.46'
$ws.Range('E48').Value = '  +2.89%  '

$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '1.60'
$ws.Range('E49').Value = '  +8.44%  '

$ws.Range('D50').Value = '2.91'
$ws.Range('E50').Value = '  +1.41%  '

$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.531.15'
$ws.Range('E51').Value = '  +1.60%  '

# Restore the default (unstyled) appearance so only cell values changed.
$ws.Range("D2:E51").Style = "Normal"
